# "Support for bar charts"
#
# 1) Turn on word-wrap for every styled header/label/data cell on "Page one"
#    (the sheet's four in-use cell formats all gain alignment/wrapText=1).
# 2) Widen columns A:F, A:K and A:P on "Page one" to make room for the
#    (future) bar-chart data ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page one")

# --- 1) wrap text on every cell block that carries an explicit style ------
# Block 1 ("Monthly Average Response Time" style cells)
$ws.Range("A1:B1").WrapText = $true
$ws.Range("A2:B2").WrapText = $true
$ws.Range("A3:A4").WrapText = $true
$ws.Range("B3:B4").WrapText = $true

# Block 2 (table with 4 columns)
$ws.Range("A6:D6").WrapText = $true
$ws.Range("A7:D7").WrapText = $true
$ws.Range("A8:D8").WrapText = $true

# Block 3
$ws.Range("A11:B11").WrapText = $true
$ws.Range("A12:B12").WrapText = $true
$ws.Range("A13:B13").WrapText = $true

# Block 4
$ws.Range("A16:B16").WrapText = $true
$ws.Range("A17:B17").WrapText = $true
$ws.Range("A18:B18").WrapText = $true

# --- 2) widen the columns that will host the new bar charts' data --------
$ws.Range("A1:F1").ColumnWidth = 29.8
$ws.Range("A1:K1").ColumnWidth = 29.8
$ws.Range("A1:P1").ColumnWidth = 29.8
